$d = $word.ActiveDocument

$replacements = @(
    @("44×80=", "39×53="),
    @("62×45=", "41×79="),
    @("95×88=", "28×48="),
    @("83×87=", "37×41="),
    @("31×54=", "77×11="),
    @("60×69=", "77×81="),
    @("64×28=", "24×52="),
    @("11×87=", "34×75="),
    @("42×91=", "56×92="),
    @("11×37=", "38×56="),
    @("50×18=", "48×57="),
    @("62×67=", "31×31="),
    @("79×18=", "51×47="),
    @("74×26=", "70×69="),
    @("61×75=", "42×26="),
    @("77×41=", "49×38="),
    @("86×59=", "91×60="),
    @("77×62=", "44×22="),
    @("42×75=", "30×78="),
    @("93×21=", "11×46="),
    @("42×77=", "71×47="),
    @("74×31=", "74×65="),
    @("17×66=", "27×25="),
    @("15×97=", "78×23="),
    @("93×30=", "85×30=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
